# Generate Report for Handoff
#
# For the rows whose Status is "Ready for handoff" (rows 7, 9, 10, 11, 12, 14
# on the per-language sheets), a fresh handoff report run:
#   - stamps a new "Latest Handoff Datetime" / "Latest HO Xliff Generate Date"
#     on the Overview sheet and on each language sheet, and
#   - marks those rows' "Priority" column as "ht" (handoff type) on each
#     language sheet.
# (Rows 8 and 13 already carry their own, different, handoff timestamp/
# priority values and are left untouched, matching the source diff.)

$wb = $excel.ActiveWorkbook

$handoffRows = @(7, 9, 10, 11, 12, 14)

# --- Overview sheet: column G = "Latest HO Xliff Generate Date" ---
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($row in $handoffRows) {
    $wsOverview.Range("G$row").Value = "2016-09-03 02:23:52"
}

# --- zh-cn sheet: column H = "Latest Handoff Datetime", column E = "Priority" ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($row in $handoffRows) {
    $wsZhCn.Range("H$row").Value = "2016-09-03 02:23:46"
    $wsZhCn.Range("E$row").Value = "ht"
}

# --- de-de sheet: column H = "Latest Handoff Datetime", column E = "Priority" ---
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($row in $handoffRows) {
    $wsDeDe.Range("H$row").Value = "2016-09-03 02:23:52"
    $wsDeDe.Range("E$row").Value = "ht"
}
